$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"25.29000000000051"
$ws.Range("H2").Value = [double]"0.0001062828338950217"
$ws.Range("I2").Value = [double]"0.0001062828338950217"
$ws.Range("L2").Value = [double]"38.87093925827847"
$ws.Range("M2").Value = "[19.90014128266175, 57.84173723389519]"
$ws.Range("N2").Value = [double]"0.0001568204925483929"
$ws.Range("O2").Value = [double]"0.0001568204925483929"
$ws.Range("P2").Value = [double]"1.628973968528041"
$ws.Range("Q2").Value = "[1.0000264903318863, 2.257921446724196]"
$ws.Range("R2").Value = [double]"4.45886434841114e-06"
$ws.Range("S2").Value = [double]"4.45886434841114e-06"
$ws.Range("T2").Value = [double]"62.4724186613745"
$ws.Range("U2").Value = "[50.503372343305486, 74.44146497944351]"
$ws.Range("V2").Value = [double]"1.06137321154165e-13"
$ws.Range("W2").Value = [double]"1.06137321154165e-13"
$ws.Range("X2").Value = [double]"18.73333333333371"
$ws.Range("Y2").Value = [double]"16.20180180180213"
$ws.Range("Z2").Value = [double]"21.2648648648653"
$ws.Range("F3").Value = [double]"25.29000000000051"
$ws.Range("H3").Value = [double]"0.01502206297587538"
$ws.Range("I3").Value = [double]"0.01502206297587538"
$ws.Range("L3").Value = [double]"29.19911776625503"
$ws.Range("M3").Value = "[3.6963544936542476, 54.70188103885582]"
$ws.Range("N3").Value = [double]"0.02577149107343613"
$ws.Range("O3").Value = [double]"0.02577149107343613"
$ws.Range("P3").Value = [double]"1.226447582482502"
$ws.Range("Q3").Value = "[0.0691842226015762, 2.383710942363427]"
$ws.Range("R3").Value = [double]"0.03828303824289958"
$ws.Range("S3").Value = [double]"0.03828303824289958"
$ws.Range("T3").Value = [double]"64.7354401903863"
$ws.Range("U3").Value = "[50.87367847519599, 78.59720190557661]"
$ws.Range("V3").Value = [double]"3.416378291376532e-12"
$ws.Range("W3").Value = [double]"3.416378291376532e-12"
$ws.Range("X3").Value = [double]"20.35351351351393"
$ws.Range("Y3").Value = [double]"15.69549549549581"
$ws.Range("Z3").Value = [double]"25.01153153153204"
$ws.Range("F4").Value = [double]"25.29000000000051"
$ws.Range("H4").Value = [double]"0.0003838091702780266"
$ws.Range("I4").Value = [double]"0.0003838091702780266"
$ws.Range("L4").Value = [double]"41.33026632083787"
$ws.Range("M4").Value = "[18.659273149263072, 64.00125949241267]"
$ws.Range("N4").Value = [double]"0.0006367954829931133"
$ws.Range("O4").Value = [double]"0.0006367954829931133"
$ws.Range("P4").Value = [double]"1.465447624197041"
$ws.Range("Q4").Value = "[0.761026448617347, 2.1698687997767347]"
$ws.Range("R4").Value = [double]"0.0001284686594120732"
$ws.Range("S4").Value = [double]"0.0001284686594120732"
$ws.Range("T4").Value = [double]"57.08785236312843"
$ws.Range("U4").Value = "[43.20487217329323, 70.97083255296363]"
$ws.Range("V4").Value = [double]"1.3439316326469e-10"
$ws.Range("W4").Value = [double]"1.3439316326469e-10"
$ws.Range("X4").Value = [double]"19.39153153153192"
$ws.Range("Y4").Value = [double]"16.55621621621655"
$ws.Range("Z4").Value = [double]"22.2268468468473"
$ws.Range("F5").Value = [double]"25.29000000000051"
$ws.Range("H5").Value = [double]"0.0001043398071879365"
$ws.Range("I5").Value = [double]"0.0001043398071879365"
$ws.Range("L5").Value = [double]"47.79382248088432"
$ws.Range("M5").Value = "[20.68850132464634, 74.8991436371223]"
$ws.Range("N5").Value = [double]"0.0009120575221266414"
$ws.Range("O5").Value = [double]"0.0009120575221266414"
$ws.Range("P5").Value = [double]"0.9119738433844242"
$ws.Range("Q5").Value = "[0.37107901213573147, 1.452868674633117]"
$ws.Range("R5").Value = [double]"0.001438956224608212"
$ws.Range("S5").Value = [double]"0.001438956224608212"
$ws.Range("T5").Value = [double]"52.49260572312033"
$ws.Range("U5").Value = "[38.257640901560116, 66.72757054468056]"
$ws.Range("V5").Value = [double]"2.380017516401267e-09"
$ws.Range("W5").Value = [double]"2.380017516401267e-09"
$ws.Range("X5").Value = [double]"21.61927927927972"
$ws.Range("Y5").Value = [double]"19.44216216216256"
$ws.Range("Z5").Value = [double]"23.79639639639688"
$ws.Range("B6").Value = [double]"0"
$ws.Range("F6").Value = [double]"25.29000000000051"
$ws.Range("H6").Value = [double]"0.000107428758387984"
$ws.Range("I6").Value = [double]"0.000107428758387984"
$ws.Range("L6").Value = [double]"49.17154188226192"
$ws.Range("M6").Value = "[21.415013285267293, 76.92807047925655]"
$ws.Range("N6").Value = [double]"0.0008681269631478639"
$ws.Range("O6").Value = [double]"0.0008681269631478639"
$ws.Range("P6").Value = [double]"0.4968685077749617"
$ws.Range("Q6").Value = "[-0.09434212172942402, 1.0880791372793475]"
$ws.Range("R6").Value = [double]"0.0974253395720226"
$ws.Range("S6").Value = [double]"0.0974253395720226"
$ws.Range("T6").Value = [double]"59.7245644022255"
$ws.Range("U6").Value = "[45.318795525218235, 74.13033327923276]"
$ws.Range("V6").Value = [double]"1.071862598678308e-10"
$ws.Range("W6").Value = [double]"1.071862598678308e-10"
$ws.Range("X6").Value = [double]"23.29009009009057"
$ws.Range("Y6").Value = [double]"20.91045045045088"
$ws.Range("Z6").Value = [double]"25.66972972973025"
$ws.Range("F7").Value = [double]"25.29000000000051"
$ws.Range("H7").Value = [double]"0.03864481252724705"
$ws.Range("I7").Value = [double]"0.03864481252724705"
$ws.Range("L7").Value = [double]"28.43266546021532"
$ws.Range("M7").Value = "[1.2727709978043151, 55.59255992262633]"
$ws.Range("N7").Value = [double]"0.04058840410917419"
$ws.Range("O7").Value = [double]"0.04058840410917419"
$ws.Range("P7").Value = [double]"-0.1761052938949232"
$ws.Range("Q7").Value = "[-1.515763422452733, 1.1635528346628865]"
$ws.Range("R7").Value = [double]"0.7923996792808596"
$ws.Range("S7").Value = [double]"0.7923996792808596"
$ws.Range("T7").Value = [double]"55.42309316855381"
$ws.Range("U7").Value = "[40.39936205482209, 70.44682428228552]"
$ws.Range("V7").Value = [double]"2.356484340992893e-09"
$ws.Range("W7").Value = [double]"2.356484340992893e-09"
$ws.Range("X7").Value = [double]"0.7088288288288425"
$ws.Range("Y7").Value = [double]"-4.68333333333343"
$ws.Range("Z7").Value = [double]"6.100990990991115"
$ws.Range("F8").Value = [double]"22.84000000000013"
$ws.Range("H8").Value = [double]"0.0008728088524112554"
$ws.Range("I8").Value = [double]"0.0008728088524112554"
$ws.Range("L8").Value = [double]"43.90451442465343"
$ws.Range("M8").Value = "[18.464721963186165, 69.3443068861207]"
$ws.Range("N8").Value = [double]"0.001139078217042178"
$ws.Range("O8").Value = [double]"0.001139078217042178"
$ws.Range("P8").Value = [double]"0.006289474781961069"
$ws.Range("Q8").Value = "[-0.7610264486173479, 0.77360539818127]"
$ws.Range("R8").Value = [double]"0.9869012674208855"
$ws.Range("S8").Value = [double]"0.9869012674208855"
$ws.Range("T8").Value = [double]"44.00013368115233"
$ws.Range("U8").Value = "[28.277200880807122, 59.72306648149755]"
$ws.Range("V8").Value = [double]"1.078561616463247e-06"
$ws.Range("W8").Value = [double]"1.078561616463247e-06"
$ws.Range("X8").Value = [double]"22.81713713713727"
$ws.Range("Y8").Value = [double]"20.02786786786798"
$ws.Range("Z8").Value = [double]"25.60640640640656"
$ws.Range("F9").Value = [double]"22.84000000000013"
$ws.Range("H9").Value = [double]"4.520937989338503e-05"
$ws.Range("I9").Value = [double]"4.520937989338503e-05"
$ws.Range("L9").Value = [double]"47.67183133600916"
$ws.Range("M9").Value = "[22.639588067797817, 72.70407460422051]"
$ws.Range("N9").Value = [double]"0.0003873465377886021"
$ws.Range("O9").Value = [double]"0.0003873465377886021"
$ws.Range("P9").Value = [double]"0.4465527095192696"
$ws.Range("Q9").Value = "[-0.09434212172942491, 0.9874475407679641]"
$ws.Range("R9").Value = [double]"0.1033020399899047"
$ws.Range("S9").Value = [double]"0.1033020399899047"
$ws.Range("T9").Value = [double]"50.37795982839627"
$ws.Range("U9").Value = "[36.69492126945646, 64.06099838733608]"
$ws.Range("V9").Value = [double]"2.47626941174417e-09"
$ws.Range("W9").Value = [double]"2.47626941174417e-09"
$ws.Range("X9").Value = [double]"21.21673673673686"
$ws.Range("Y9").Value = [double]"19.25053053053064"
$ws.Range("Z9").Value = [double]"23.18294294294308"
$ws.Range("B10").Value = [double]"1"
$ws.Range("F10").Value = [double]"22.84000000000013"
$ws.Range("H10").Value = [double]"0.002463556132360045"
$ws.Range("I10").Value = [double]"0.002463556132360045"
$ws.Range("L10").Value = [double]"41.13295238875114"
$ws.Range("M10").Value = "[10.065837503256517, 72.20006727424577]"
$ws.Range("N10").Value = [double]"0.01060572768814994"
$ws.Range("O10").Value = [double]"0.01060572768814994"
$ws.Range("P10").Value = [double]"0.8365001460008861"
$ws.Range("Q10").Value = "[0.16981581911296217, 1.50318447288881]"
$ws.Range("R10").Value = [double]"0.01508310376789823"
$ws.Range("S10").Value = [double]"0.01508310376789823"
$ws.Range("T10").Value = [double]"66.01511806945804"
$ws.Range("U10").Value = "[50.320131970425024, 81.71010416849106]"
$ws.Range("V10").Value = [double]"7.170108951015663e-11"
$ws.Range("W10").Value = [double]"7.170108951015663e-11"
$ws.Range("X10").Value = [double]"19.79923923923935"
$ws.Range("Y10").Value = [double]"17.37577577577587"
$ws.Range("Z10").Value = [double]"22.22270270270283"

Write-Output "Applied 164 cell updates"
